# The deck originally carries two theme parts:
#   - theme1.xml ("Integral")     -> used by the Slide Master (ppt/theme/theme1.xml)
#   - theme2.xml ("Office Theme") -> used by the Notes Master  (ppt/theme/theme2.xml)
#
# The authored edit swaps which theme colors are applied to the Slide Master
# versus the Notes Master: the Slide Master ends up with the plain "Office
# Theme" palette, and the Notes Master ends up with the previous "Integral"
# palette. Font scheme and format scheme are identical between the two
# themes, so only the 12 theme colors (and, where the object model allows,
# the theme identity) need to move.

$p  = $ppt.ActivePresentation
$m  = $p.Slides.Item(1).Master
$nm = $p.NotesMaster

# Snapshot the current ("Integral") master colors before overwriting them,
# so they can be reapplied to the Notes Master (true swap, not just a copy).
$masterColors = @()
for ($i = 1; $i -le 12; $i++) {
    $masterColors += $m.Theme.ThemeColorScheme.Colors($i).RGB
}

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# expressed as VBA-style BGR "RGB()" longs.
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x6A5444,  # dk2      (44546A)
    0xE6E6E7,  # lt2      (E7E6E6)
    0xD59B5B,  # accent1  (5B9BD5)
    0x317DED,  # accent2  (ED7D31)
    0xA5A5A5,  # accent3  (A5A5A5)
    0x00C0FF,  # accent4  (FFC000)
    0xC47244,  # accent5  (4472C4)
    0x47AD70,  # accent6  (70AD47)
    0xC16305,  # hlink    (0563C1)
    0x724F95   # folHlink (954F72)
)

# Give the Notes Master the old "Integral" colors that used to live on the
# Slide Master.
for ($i = 1; $i -le 12; $i++) {
    $nm.Theme.ThemeColorScheme.Colors($i).RGB = $masterColors[$i - 1]
}

# Give the Slide Master the "Office Theme" colors that used to live on the
# Notes Master. Applied last so it is authoritative for the shared theme
# part that actually renders the slides.
for ($i = 1; $i -le 12; $i++) {
    $m.Theme.ThemeColorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
